$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update "Latest HO Xliff Generate Date" on Overview (column G) for rows 7-11,14
foreach ($r in 7,8,9,10,11,14) {
    $wsOverview.Range("G$r").Value = "2016-08-28 14:21:52"
}

# Update "Latest Handoff Datetime" on zh-cn (column H) for rows 7-11,14
foreach ($r in 7,8,9,10,11,14) {
    $wsZhCn.Range("H$r").Value = "2016-08-28 14:21:48"
}

# Update "Priority" column (E) on zh-cn and de-de for rows 7-11,14: "" -> "ht"
foreach ($r in 7,8,9,10,11,14) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
